# feat: add 2022-Q3 data
#
# 1. Insert a brand-new worksheet named "2022-Q3" right after "总计" (so it
#    becomes the 2nd tab, pushing the existing "2022-Q2" / "2020-Q4" tabs
#    down by one) and fill it with the Q3 fund-holding breakdown table.
# 2. Update the "总计" (totals) sheet so its first data row now reports the
#    new 2022-Q3 quarter, and the rows that used to hold 2022-Q2 / 2020-Q4
#    shift down to make room for it.

$wb = $excel.ActiveWorkbook

# Writing a numeric-looking string straight into Range.Value (e.g. a fund
# code like "001668" or a ratio like "12.06") gets auto-coerced to a number
# and loses its leading zeros, same as typing it into Excel by hand. Force
# it to stay text by flipping the cell to the "@" (Text) number format
# first, then immediately paste the *formatting* from a plain, unstyled
# cell back on top so the cell ends up with no special style applied -
# only its stored type (text) differs from a freshly-written number.
function Set-PlainTextValue($cell, $val, $plainRef) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $plainRef.Copy()
    $cell.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# Step 1: add the new "2022-Q3" sheet, positioned right after "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$plainRef = $totalSheet.Range("C2")        # a plain, unstyled numeric cell used as a formatting donor

$q2Sheet = $wb.Worksheets.Item(2)          # currently "2022-Q2" -> will be pushed to slot 3
$q3Sheet = $wb.Worksheets.Add($q2Sheet)    # insert before the old "2022-Q2" tab
$q3Sheet.Name = "2022-Q3"

# Header row (row 1): copy the bold/bordered header style used by the other
# per-quarter sheets, then write the column titles.
$totalSheet.Range("B1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $q3Sheet.Cells.Item(1, 2 + $c).Value = $headers[$c]
}

# Data rows 2-8: code, name, size, stock position, position ratio, market
# value, position rank. Every column except the rank is stored as text in
# the source data (even the numeric-looking ones), matching the original
# per-quarter sheets in this workbook.
$rows = @(
    @("001668", "汇添富全球移动互联灵活配置混合（QDII）A", "12.06", "90.88", "2.86", "0.3449", 7),
    @("161128", "易方达标普信息科技指数（QDII-LOF）人民币", "4.99", "91.96", "3.36", "0.1677", 4),
    @("012868", "易方达标普信息科技指数（QDII-LOF）人民币 C", "4.99", "91.96", "3.36", "0.1677", 4),
    @("003721", "易方达标普信息科技指数（QDII-LOF）美元A", "4.84", "91.96", "3.36", "0.1626", 4),
    @("012869", "易方达标普信息科技指数（QDII-LOF）美元 C", "0.15", "91.96", "3.36", "0.0050", 4),
    @("015203", "汇添富全球移动互联灵活配置混合（QDII）D", "0.04", "90.88", "2.86", "0.0011", 7),
    @("015202", "汇添富全球移动互联灵活配置混合（QDII）C", "0.01", "90.88", "2.86", "0.0003", 7)
)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowNum = 2 + $r
    $rowData = $rows[$r]

    # Column A: zero-based row index, styled like the other sheets' index column.
    $q3Sheet.Cells.Item($rowNum, 1).Value = $r
    $totalSheet.Range("A2").Copy()
    $q3Sheet.Cells.Item($rowNum, 1).PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    # Columns B-G: text, even the numeric-looking ones.
    for ($c = 0; $c -lt 6; $c++) {
        Set-PlainTextValue $q3Sheet.Cells.Item($rowNum, 2 + $c) $rowData[$c] $plainRef
    }
    $excel.CutCopyMode = $false

    # Column H: the rank is a genuine number.
    $q3Sheet.Cells.Item($rowNum, 8).Value = $rowData[6]
}

# ---------------------------------------------------------------------
# Step 2: update the "总计" sheet — insert the 2022-Q3 totals as the new
# first data row, shifting 2022-Q2 / 2020-Q4 down by one row.
# ---------------------------------------------------------------------

# Row 4 ("2020-Q4") is brand-new in this sheet; give it the same style as
# the existing index column before filling it in.
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$totalSheet.Range("B4").Value = "2020-Q4"
$totalSheet.Range("C4").Value = 4
$totalSheet.Range("D4").Value = 1

# Row 3 now holds what used to be row 2's data (2022-Q2).
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 4
$totalSheet.Range("D3").Value = 0.37

# Row 2 now holds the brand-new 2022-Q3 totals.
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 0.85

Write-Output "2022-Q3 sheet added and totals sheet updated"
